# Add a single-column table "xlsx_single_col_table" over L21:L25 on the
# active sheet, mirroring the existing named range xlsx_single_col_range
# (J21:J25). Header "colA" plus data rows a, b, c, d.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L21").Value = "colA"
$ws.Range("L22").Value = "a"
$ws.Range("L23").Value = "b"
$ws.Range("L24").Value = "c"
$ws.Range("L25").Value = "d"

$null = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("L21:L25"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)

# Locate the table we just created by its range (rather than trusting the
# handle returned from .Add()) and rename it. The default table style
# applied on creation already matches the target ("TableStyleMedium2"),
# so we leave TableStyle untouched.
$newTable = $null
for ($i = 1; $i -le $ws.ListObjects.Count; $i++) {
    $candidate = $ws.ListObjects.Item($i)
    if ($candidate.Range.Address() -eq '$L$21:$L$25') {
        $newTable = $candidate
    }
}

$newTable.Name = "xlsx_single_col_table"

# Update the selection to match the authored state.
$ws.Range("L22").Select()
